# LOM3046.docx - Bibliografia paragraph: turn the single run-on paragraph
# listing references 1-19 into one run per reference separated by manual
# line breaks (<w:br/>), mirroring the target diff.
#
# We locate each "boundary" between two consecutive numbered references
# (a short, unique slice of text spanning the end of reference N and the
# start of reference N+1) and replace it with the same text plus a Word
# manual line break ("^l" in Find/Replace syntax, which yields <w:br/>).
# Because each boundary string is unique within the document, this is
# equivalent to precisely splitting the big run at 18 points, turning the
# single paragraph into 19 "lines" joined by line breaks - without
# touching any other paragraph.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$boundaries = @(
    @("Campus, Rio de Janeiro, 1984. ", "2. Shackelford, J.F."),
    @("on. Prentice Hall Inc., 1996. ", "3. Padilha, A.F. Téc"),
    @(", Ed. Hemus, São Paulo, 1985. ", "4. Guy, A.G. Ciência"),
    @(" e Científicos Editora, 1982. ", "5. Reed-Hill, R.E. P"),
    @("ca, Ed. Guanabara Dois, 1982. ", "6. Nondestructive Ch"),
    @("ries. Plenum Press, New York. ", "7. Yacobi, B.G. Holt"),
    @("Plenum Press, New York, 1994. ", "8. Lowell, S.; Shiel"),
    @(" and Density, Springer, 2010. ", "9. Murphy, D. B. Fun"),
    @("ic Imaging, Wiley-Liss, 2001. ", "10. Wu, Q.; Merchant"),
    @("essing, Academic Press, 2008. ", "11. Cullity, B. D.; "),
    @("raction, Prentice Hall, 2001. ", "12. Goldstein, J.; e"),
    @("icroanalysis, Springer, 2003. ", "13. Hatakeyama, T.; "),
    @("al Analysis, NY: Wiley, 1999. ", "14. Haines, P. J. Pr"),
    @("l Society of Chemistry, 2002. ", "15. Schramm, G. Reol"),
    @("etria. Editora Artliber, 2006.", "16. Azevedo, A. D.; "),
    @("is. São Paulo: ARTLIBER, 2009.", "17. Brown, M.E. Hand"),
    @("erdam: Elsevier Science, 1998.", "18. Muller, A. Solid"),
    @("Porto Alegre: Ed. UFRGS, 2002.", "19. Speyer, R. Therm")
)

$i = 0
foreach ($b in $boundaries) {
    $i = $i + 1
    $left = $b[0]
    $right = $b[1]
    $needle = $left + $right
    $replacement = $left + "^l" + $right
    $ok = $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
    if (-not $ok) {
        Write-Output "WARNING: boundary $i not found/replaced: $needle"
    }
}

Write-Output "Bibliografia references split into $($i) line breaks."
